$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Write new cell values in the same order the shared-string table picks them up:
# column C first, then B, then A, matching the original commit's string ordering.

# Row 34
$ws.Range("C34").Value = "Verify that anyone can see the public watchlists of a user on user's profile page"
$ws.Range("B34").Value = "OPQA-321"
$ws.Range("A34").Value = "TestCase_E33"
$ws.Range("D34").Value = "Y"
$ws.Range("E34").Value = "PASS"

# Row 35
$ws.Range("C35").Value = "Verify that no one can see the private watchlists of a user on user's profile page"
$ws.Range("B35").Value = "OPQA-329"
$ws.Range("A35").Value = "TestCase_E34"
$ws.Range("D35").Value = "Y"
$ws.Range("E35").Value = "PASS"

# Apply the same styling pattern as surrounding rows:
# A/B/D/E use the bordered style (row 33), C uses the wrap+fill style (row 32)
$ws.Range("A33:B33").Copy() | Out-Null
$ws.Range("A34:B34").PasteSpecial(-4122) | Out-Null
$ws.Range("A35:B35").PasteSpecial(-4122) | Out-Null

$ws.Range("D33:E33").Copy() | Out-Null
$ws.Range("D34:E34").PasteSpecial(-4122) | Out-Null
$ws.Range("D35:E35").PasteSpecial(-4122) | Out-Null

$ws.Range("C32").Copy() | Out-Null
$ws.Range("C34").PasteSpecial(-4122) | Out-Null
$ws.Range("C35").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Update the sheet view to match the new extent: top-left at C1, selection E2:E35
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("E2:E35").Select() | Out-Null
